$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.826.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "'2.248.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'317.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "'101.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  -1.47%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.556"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").Value = "'36.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("D11").Value = "'0.0831"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "'7.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").Value = "'2.588.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "'0.857"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'14.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.243.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").Value = "'43.719.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").Value = "'13.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.78%  "
$ws.Range("D20").Value = "'0.0₃0985"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").Value = "'6.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("D22").Value = "'65.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'3.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D24").Value = "'235.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("E25").Value = "  -2.94%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'10.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'2.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").Value = "'37.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("D30").Value = "'6.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.46%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'20.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'158.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("D33").Value = "'0.0850"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("E34").Value = "  -0.73%  "
$ws.Range("E35").Value = "  +10.97%  "
$ws.Range("D36").Value = "'3.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.82%  "
$ws.Range("D37").Value = "'1.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").Value = "'0.119"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("E40").Value = "  -6.04%  "
$ws.Range("D41").Value = "'15.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +16.25%  "
$ws.Range("D42").Value = "'0.0316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").Value = "'1.789.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.198"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("B46").Value = "ordi"
$ws.Range("C46").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D46").Value = "'75.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").Value = "'82.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.19%  "
$ws.Range("D48").Value = "'5.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").Value = "'58.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").Value = "'103.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "'1.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.85%  "
